$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column AD: new day "15-10-2020" added to the COVID19 time series ---

# AD1: header label, formatted like the rest of the date header row (N1:AC1)
# i.e. bold font + thin box border + centered/top alignment.
$headerCell = $ws.Cells.Item(1, 30)
$headerCell.Value = "15-10-2020"
$headerCell.Font.Bold = $true
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = -4108   # xlCenter
$headerCell.VerticalAlignment = -4160     # xlTop

# AD2:AD36: cumulative case counts per state/UT for 15-10-2020
$ws.Cells.Item(2, 30).Value = 3796
$ws.Cells.Item(3, 30).Value = 719477
$ws.Cells.Item(4, 30).Value = 9694
$ws.Cells.Item(5, 30).Value = 168072
$ws.Cells.Item(6, 30).Value = 188380
$ws.Cells.Item(7, 30).Value = 12119
$ws.Cells.Item(8, 30).Value = 121548
$ws.Cells.Item(9, 30).Value = 3090
$ws.Cells.Item(10, 30).Value = 289747
$ws.Cells.Item(11, 30).Value = 34731
$ws.Cells.Item(12, 30).Value = 136404
$ws.Cells.Item(13, 30).Value = 133706
$ws.Cells.Item(14, 30).Value = 15233
$ws.Cells.Item(15, 30).Value = 74318
$ws.Cells.Item(16, 30).Value = 86367
$ws.Cells.Item(17, 30).Value = 611167
$ws.Cells.Item(18, 30).Value = 215149
$ws.Cells.Item(19, 30).Value = 4261
$ws.Cells.Item(20, 30).Value = 138158
$ws.Cells.Item(21, 30).Value = 1316769
$ws.Cells.Item(22, 30).Value = 10915
$ws.Cells.Item(23, 30).Value = 5582
$ws.Cells.Item(24, 30).Value = 2108
$ws.Cells.Item(25, 30).Value = 5916
$ws.Cells.Item(26, 30).Value = 235763
$ws.Cells.Item(27, 30).Value = 27152
$ws.Cells.Item(28, 30).Value = 114075
$ws.Cells.Item(29, 30).Value = 141835
$ws.Cells.Item(30, 30).Value = 3075
$ws.Cells.Item(31, 30).Value = 617403
$ws.Cells.Item(32, 30).Value = 193218
$ws.Cells.Item(33, 30).Value = 25394
$ws.Cells.Item(34, 30).Value = 49129
$ws.Cells.Item(35, 30).Value = 401306
$ws.Cells.Item(36, 30).Value = 268384

